$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.280.69'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '2.021.78'
$ws.Range("E3").Value = '  +2.61%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.92'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.45%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +2.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0800'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.92'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.87%  '
$ws.Range("D13").Value = '2.322.82'
$ws.Range("E13").Value = '  +2.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.832'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.48'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.38'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").Value = '2.030.42'
$ws.Range("E17").Value = '  +3.37%  '
$ws.Range("D18").Value = '37.268.62'
$ws.Range("E18").Value = '  +2.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.84'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '0.0₃0853'
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.33'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E24").Value = '  +4.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.10'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.137'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -5.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.78'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.36'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.121'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("E33").Value = '  +8.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.54'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  +8.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.55'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.11%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +2.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.35'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0969'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  +3.39%  '
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("E44").Value = '  +3.21%  '
$ws.Range("D45").Value = '1.392.21'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.64'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.43'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("E49").Value = '  +11.67%  '
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").Value = '2.214.38'
$ws.Range("E51").Value = '  +2.94%  '
